$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = 20
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = 0.486
$ws.Range("E5").Value = 0.092
$ws.Range("F5").Value = 0.169
$ws.Range("G5").Formula = "=D5"
$ws.Range("H5").Value = 55.5
$ws.Range("I5").Value = 14
$ws.Range("J5").Formula = "=1000/I5"
$ws.Range("L5").Value = "Simon's account"

$ws.Range("G6").Select()
